# Edits the titles of slides 6, 7, 8 and 9 ("Sigmoid 이용", "Relu 이용",
# "학습률 조정", "입력노드개수, 학습률 조정, Dropout") to append the
# hidden-layer / input-node / learning-rate experiment details, and turns
# on "shrink text on overflow" for each of those title placeholders.

$p = $ppt.ActivePresentation
$vt = [char]0x0B   # soft line-break marker recognised by InsertAfter

# ---------------------------------------------------------------------
# Slide 6 : "Sigmoid 이용"  ->  "Sigmoid 이용, hidden layer 1층, <br>input node 500개,  learning rate 0.1"
# ---------------------------------------------------------------------
$s6 = $p.Slides.Item(6)
$t6 = $s6.Shapes.Item(1).TextFrame
$t6.AutoSize = 2
$tr6 = $t6.TextRange
$tr6.InsertAfter(", hidden layer 1")
$tr6.InsertAfter("층")
$tr6.InsertAfter(", ")
$tr6.InsertAfter("${vt}input node 500")
$tr6.InsertAfter("개")
$tr6.InsertAfter(",  learning rate 0.1")

# ---------------------------------------------------------------------
# Slide 7 : "Relu 이용"  ->  "Relu 이용, hidden layer 1층, <br>input node 500개,  learning rate 0.1"
# ---------------------------------------------------------------------
$s7 = $p.Slides.Item(7)
$t7 = $s7.Shapes.Item(1).TextFrame
$t7.AutoSize = 2
$tr7 = $t7.TextRange
$tr7.InsertAfter(", hidden layer 1")
$tr7.InsertAfter("층")
$tr7.InsertAfter(", ")
$tr7.InsertAfter("${vt}input node 500")
$tr7.InsertAfter("개")
$tr7.InsertAfter(",  learning rate 0.1")

# ---------------------------------------------------------------------
# Slide 8 : "학습률 조정"  ->  "Relu이용, hidden layer 1층, <br>input node 500개,  learning rate 0.01"
# ---------------------------------------------------------------------
$s8 = $p.Slides.Item(8)
$t8 = $s8.Shapes.Item(1).TextFrame
$t8.AutoSize = 2
$tr8 = $t8.TextRange
$tr8.Text = "Relu"
$tr8.LanguageID = "en-US"
$tr8.InsertAfter("이용")
$tr8.InsertAfter(", hidden layer 1")
$tr8.InsertAfter("층")
$tr8.InsertAfter(", ")
$tr8.InsertAfter("${vt}input node 500")
$tr8.InsertAfter("개")
$tr8.InsertAfter(",  learning rate ")
$tr8.InsertAfter("0.01")

# ---------------------------------------------------------------------
# Slide 9 : "입력노드개수, 학습률 조정, Dropout" ->
#           "Relu이용, hidden layer 1층, dropout 0.7,  <br>input node 700개,  learning rate 0.001,"
# ---------------------------------------------------------------------
$s9 = $p.Slides.Item(9)
$t9 = $s9.Shapes.Item(1).TextFrame
$t9.AutoSize = 2
$tr9 = $t9.TextRange
$tr9.Text = "Relu"
$tr9.LanguageID = "en-US"
$tr9.InsertAfter("이용")
$tr9.InsertAfter(", hidden layer 1")
$tr9.InsertAfter("층")
$tr9.InsertAfter(", dropout 0.7,  ")
$tr9.InsertAfter("${vt}input node ")
$tr9.InsertAfter("700")
$tr9.InsertAfter("개")
$tr9.InsertAfter(",  learning rate ")
$tr9.InsertAfter("0.001,")
